# Updated cryptos list on Tue Jul 25 19:00:45 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds values like "1.001" or "0.7002" that Excel's COM
# layer would otherwise auto-convert into numbers/dates. Force the whole
# price column to Text format first so the new values are written verbatim,
# matching the existing inlineStr text cells in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.285.62"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "1.861.52"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "0.7013"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "237.99"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "0.08218"
$ws.Range("E8").Value = "  +9.55%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "23.31"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "1.879.40"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("D13").Value = "0.7173"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "5.180"
$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").Value = "89.37"

$ws.Range("D16").Value = "29.300.29"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "5.784"
$ws.Range("E17").Value = "  -0.18%  "

# Rows 18 and 19 swap places: ShibaInu <-> Avalanche
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "13.40"
$ws.Range("E18").Value = "  +2.50%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007862"
$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("D20").Value = "237.63"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "2.107.26"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "7.473"
$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "162.30"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").Value = "8.992"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "0.1439"
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").Value = "18.13"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "1.980"
$ws.Range("E29").Value = "  +2.05%  "

$ws.Range("D30").Value = "1.430"
$ws.Range("E30").Value = "  +3.29%  "

$ws.Range("D31").Value = "4.436"
$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("D33").Value = "4.061"
$ws.Range("E33").Value = "  +1.29%  "

$ws.Range("D34").Value = "0.05227"
$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("D35").Value = "1.172"
$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("D36").Value = "0.7063"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("D38").Value = "2.663"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("E39").Value = "  -0.60%  "

$ws.Range("D40").Value = "2.732"
$ws.Range("E40").Value = "  +1.96%  "

$ws.Range("D41").Value = "1.144.23"
$ws.Range("E41").Value = "  +6.16%  "

$ws.Range("D42").Value = "0.9195"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "5.981"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "0.4284"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").Value = "70.95"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "102.72"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("D49").Value = "2.004.19"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").Value = "9.184"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").Value = "6.975"
$ws.Range("E51").Value = "  -1.36%  "
